$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.208.06'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").Value = '2.380.34'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("D9").Value = '2.378.76'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.339'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000167'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.794.33'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '60.043.67'
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").Value = '2.376.62'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.10'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +12.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '559.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.64%  '
$ws.Range("D29").Value = '2.495.86'
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("D30").Value = '0.0₃0931'
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.64%  '
$ws.Range("E32").Value = '  -2.52%  '
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '151.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.47%  '
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.70%  '
$ws.Range("D46").Value = '0.0₆0288'
$ws.Range("E46").Value = '  +1.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.29%  '
